$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "department" value for the course row from
# "EDISON SCHOOL OF TECH SCIENCES" to "Packages"
$ws.Range("C2").Value = "Packages"
